$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the new shared-strings in the same order the original edit created
# them (a, b, c, e, d) so the resulting string table matches byte-for-byte.
$ws.Cells.Item(2, 9).Value = "a"
$ws.Cells.Item(3, 9).Value = "b"
$ws.Cells.Item(4, 9).Value = "c"
$ws.Cells.Item(6, 9).Value = "e"
$ws.Cells.Item(5, 9).Value = "d"

# The cyclic letters (a,b,c,d,e) assigned to column I starting at row 2.
$letters = @("a", "b", "c", "d", "e")

for ($row = 2; $row -le 91; $row++) {
    $letter = $letters[($row - 2) % 5]
    $cell = $ws.Cells.Item($row, 9)
    $cell.Value = $letter
    $cell.Style = "Hipervínculo"
}

# Update the view: move the active selection to J6 (this also clears the
# previous scrolled-down topLeftCell position recorded in the sheet view).
$ws.Range("J6").Select()
